$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Clase: Arbitro" intro paragraph.
# Word's grammar checker flags "de acuerdo a" (gramStart/gramEnd) splitting
# the single run into three runs.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Es una de las clases principales", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2FA968A1" w14:textId="2E28AD55" w:rsidR="00013013" w:rsidRDefault="00013013"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Es una de las clases principales del programa ya que cumple un papel fundamental para el torneo, evaluando los lanzamientos </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>de acuerdo a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> un criterio.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: Arbitro method "evaluarSegun" renamed to "cambiarCriterioDeEvaluacion"
# and its description expanded to mention the circular queue.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("evaluarSegun", $true, $false, $false, $false, $false, $true, 1, $false, "cambiarCriterioDeEvaluacion", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("Strategy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$period = $d.Range($rng2.End, $rng2.End + 1)
$period.Text = ", con una cola circular, que cambia al siguiente criterio de evaluación."

# ---------------------------------------------------------------------------
# Edit 3: "Calcular" interface paragraph — grammar mark around "dos método".
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("la cual tiene", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(1)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3E124384" w14:textId="3A33DD05" w:rsidR="005033D0" w:rsidRPr="00DF0E2C" w:rsidRDefault="005033D0" w:rsidP="005033D0"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00DF0E2C"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Esta es una </w:t></w:r><w:r w:rsidR="00314325"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>interface</w:t></w:r><w:r w:rsidRPr="00DF0E2C"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> la cual tiene </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00926F81"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>dos</w:t></w:r><w:r w:rsidRPr="00DF0E2C"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> método</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00DF0E2C"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3) | Out-Null

Write-Output "All edits applied"
